$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.96787706683382
$ws.Range("C2").Value = 10.01888404129415
$ws.Range("D2").Value = 8.275664005662895
$ws.Range("F2").Value = 40.86700472923994
$ws.Range("G2").Value = 3.710722115901969
$ws.Range("J2").Value = 11.29545613467541
$ws.Range("K2").Value = 9.978001582170828
$ws.Range("L2").Value = 11.28768206590767
$ws.Range("M2").Value = 15.73606501363036
$ws.Range("O2").Value = 31.42334705877943
$ws.Range("B3").Value = 13.80228454092903
$ws.Range("C3").Value = 10.01754690130858
$ws.Range("D3").Value = 8.267710291138874
$ws.Range("F3").Value = 40.9612298958267
$ws.Range("G3").Value = 3.712618085478408
$ws.Range("J3").Value = 11.31700655144978
$ws.Range("K3").Value = 9.851033364114198
$ws.Range("L3").Value = 11.2961030625683
$ws.Range("M3").Value = 15.71755528567898
$ws.Range("O3").Value = 31.51116197091622
$ws.Range("B4").Value = 13.70213551727004
$ws.Range("C4").Value = 10.01696648558508
$ws.Range("D4").Value = 8.263874005294101
$ws.Range("F4").Value = 41.02636603321042
$ws.Range("G4").Value = 3.713844669992624
$ws.Range("J4").Value = 11.33099274025688
$ws.Range("K4").Value = 9.773831945026885
$ws.Range("L4").Value = 11.30246674823768
$ws.Range("M4").Value = 15.70819261401253
$ws.Range("O4").Value = 31.57020805053934
$ws.Range("B5").Value = 13.66175279381113
$ws.Range("C5").Value = 10.01679083377262
$ws.Range("D5").Value = 8.262575674668749
$ws.Range("F5").Value = 41.05473893987782
$ws.Range("G5").Value = 3.714360265701639
$ws.Range("J5").Value = 11.33688235894935
$ws.Range("K5").Value = 9.742597136911249
$ws.Range("L5").Value = 11.30536072068499
$ws.Range("M5").Value = 15.70488419558277
$ws.Range("O5").Value = 31.59555839726432
$ws.Range("B6").Value = 13.65507448078392
$ws.Range("C6").Value = 10.0167653543829
$ws.Range("D6").Value = 8.26237613770636
$ws.Range("F6").Value = 41.05956068038253
$ws.Range("G6").Value = 3.714446832822881
$ws.Range("J6").Value = 11.33787182386611
$ws.Range("K6").Value = 9.737425248630704
$ws.Range("L6").Value = 11.30585944471662
$ws.Range("M6").Value = 15.70436554906077
$ws.Range("O6").Value = 31.59984560852728
$ws.Range("B7").Value = 13.701589106237
$ws.Range("C7").Value = 10.01696386971802
$ws.Range("D7").Value = 8.26385542055921
$ws.Range("F7").Value = 41.02674127523918
$ws.Range("G7").Value = 3.713851559650607
$ws.Range("J7").Value = 11.33107139921041
$ws.Range("K7").Value = 9.773409742954934
$ws.Range("L7").Value = 11.30250455887572
$ws.Range("M7").Value = 15.70814593856298
$ws.Range("O7").Value = 31.5705447177399
$ws.Range("B8").Value = 13.91049148073714
$ws.Range("C8").Value = 10.01837328280367
$ws.Range("D8").Value = 8.27270520149054
$ws.Range("F8").Value = 40.8979813902405
$ws.Range("G8").Value = 3.711362913706171
$ws.Range("J8").Value = 11.30273049496842
$ws.Range("K8").Value = 9.934086706965624
$ws.Range("L8").Value = 11.29033840820968
$ws.Range("M8").Value = 15.72926939354787
$ws.Range("O8").Value = 31.45256103191173
$ws.Range("B9").Value = 14.33024036308497
$ws.Range("C9").Value = 10.02303030510562
$ws.Range("D9").Value = 8.298294176532755
$ws.Range("F9").Value = 40.70332444864407
$ws.Range("G9").Value = 3.706975993407841
$ws.Range("J9").Value = 11.25311524584119
$ws.Range("K9").Value = 10.25364566887857
$ws.Range("L9").Value = 11.27591688379726
$ws.Range("M9").Value = 15.7864211711472
$ws.Range("O9").Value = 31.2619101836838
$ws.Range("B10").Value = 14.64207485287662
$ws.Range("C10").Value = 10.02758547885549
$ws.Range("D10").Value = 8.322007402047674
$ws.Range("F10").Value = 40.59565170293505
$ws.Range("G10").Value = 3.704050551950063
$ws.Range("J10").Value = 11.22026541721034
$ws.Range("K10").Value = 10.48909865411634
$ws.Range("L10").Value = 11.27103144471979
$ws.Range("M10").Value = 15.83777354223189
$ws.Range("O10").Value = 31.14670101993051
$ws.Range("B11").Value = 14.78411374592527
$ws.Range("C11").Value = 10.0298995905776
$ws.Range("D11").Value = 8.333836801766163
$ws.Range("F11").Value = 40.55435551511913
$ws.Range("G11").Value = 3.702783660296766
$ws.Range("J11").Value = 11.20609673068249
$ws.Range("K11").Value = 10.59593066668795
$ws.Range("L11").Value = 11.27003899253023
$ws.Range("M11").Value = 15.86311516227385
$ws.Range("O11").Value = 31.0996958412807
$ws.Range("B12").Value = 14.83787970835728
$ws.Range("C12").Value = 10.03081029581983
$ws.Range("D12").Value = 8.338463800239465
$ws.Range("F12").Value = 40.53982338918863
$ws.Range("G12").Value = 3.702313061099749
$ws.Range("J12").Value = 11.20084233842701
$ws.Range("K12").Value = 10.63631054249683
$ws.Range("L12").Value = 11.26983912741171
$ws.Range("M12").Value = 15.87299134493993
$ws.Range("O12").Value = 31.08267389489446
$ws.Range("B13").Value = 14.82630192984431
$ws.Range("C13").Value = 10.0306126353535
$ws.Range("D13").Value = 8.337460774320885
$ws.Range("F13").Value = 40.54290394999028
$ws.Range("G13").Value = 3.702414007027128
$ws.Range("J13").Value = 11.20196903686975
$ws.Range("K13").Value = 10.62761790505798
$ws.Range("L13").Value = 11.26987436129323
$ws.Range("M13").Value = 15.87085196795174
$ws.Range("O13").Value = 31.08630526654369
$ws.Range("B14").Value = 14.78853779515099
$ws.Range("C14").Value = 10.02997382762605
$ws.Range("D14").Value = 8.334214526500681
$ws.Range("F14").Value = 40.55313778200509
$ws.Range("G14").Value = 3.702744760788092
$ws.Range("J14").Value = 11.20566222694844
$ws.Range("K14").Value = 10.59925444701413
$ws.Range("L14").Value = 11.27001902848818
$ws.Range("M14").Value = 15.86392210256577
$ws.Range("O14").Value = 31.0982798410078
$ws.Range("B15").Value = 14.76540203331236
$ws.Range("C15").Value = 10.02958700794635
$ws.Range("D15").Value = 8.332245236738403
$ws.Range("F15").Value = 40.55955033220688
$ws.Range("G15").Value = 3.702948546675843
$ws.Range("J15").Value = 11.20793885255492
$ws.Range("K15").Value = 10.58187019998926
$ws.Range("L15").Value = 11.27013052803121
$ws.Range("M15").Value = 15.85971365286363
$ws.Range("O15").Value = 31.10571594101413
$ws.Range("B16").Value = 14.63279189218244
$ws.Range("C16").Value = 10.02743908140948
$ws.Range("D16").Value = 8.321255091188515
$ws.Range("F16").Value = 40.59850521555281
$ws.Range("G16").Value = 3.704134628446826
$ws.Range("J16").Value = 11.22120692850722
$ws.Range("K16").Value = 10.48210832388084
$ws.Range("L16").Value = 11.271120976825
$ws.Range("M16").Value = 15.83615684947692
$ws.Range("O16").Value = 31.14988173202227
$ws.Range("B17").Value = 14.55145284037964
$ws.Range("C17").Value = 10.0261830880603
$ws.Range("D17").Value = 8.314778240617141
$ws.Range("F17").Value = 40.62437165748176
$ws.Range("G17").Value = 3.704878587378162
$ws.Range("J17").Value = 11.2295446084038
$ws.Range("K17").Value = 10.42081136384886
$ws.Range("L17").Value = 11.27204297272557
$ws.Range("M17").Value = 15.82220943582291
$ws.Range("O17").Value = 31.17836071063084
$ws.Range("B18").Value = 14.50468886068101
$ws.Range("C18").Value = 10.02548346982889
$ws.Range("D18").Value = 8.311151168949445
$ws.Range("F18").Value = 40.63997263859881
$ws.Range("G18").Value = 3.705312510680806
$ws.Range("J18").Value = 11.23441318153233
$ws.Range("K18").Value = 10.38553098418157
$ws.Range("L18").Value = 11.27268910728649
$ws.Range("M18").Value = 15.81437404542496
$ws.Range("O18").Value = 31.19524967679157
$ws.Range("B19").Value = 14.4888602948294
$ws.Range("C19").Value = 10.02525051818417
$ws.Range("D19").Value = 8.309940048970716
$ws.Range("F19").Value = 40.64537905954509
$ws.Range("G19").Value = 3.705460464580466
$ws.Range("J19").Value = 11.23607414128979
$ws.Range("K19").Value = 10.37358262417094
$ws.Range("L19").Value = 11.27292779776239
$ws.Range("M19").Value = 15.81175334809921
$ws.Range("O19").Value = 31.20105531952095
$ws.Range("B20").Value = 14.5601097891934
$ws.Range("C20").Value = 10.02631443422293
$ws.Range("D20").Value = 8.315457561790078
$ws.Range("F20").Value = 40.62154326637202
$ws.Range("G20").Value = 3.704798769182825
$ws.Range("J20").Value = 11.22864950026266
$ws.Range("K20").Value = 10.42733927135878
$ws.Range("L20").Value = 11.27193284350082
$ws.Range("M20").Value = 15.82367486511762
$ws.Range("O20").Value = 31.17527642939875
$ws.Range("B21").Value = 14.79963098282074
$ws.Range("C21").Value = 10.03016053040625
$ws.Range("D21").Value = 8.335164045780271
$ws.Range("F21").Value = 40.5501018401649
$ws.Range("G21").Value = 3.70264736254589
$ws.Range("J21").Value = 11.20457443898995
$ws.Range("K21").Value = 10.60758779658875
$ws.Range("L21").Value = 11.26997176834754
$ws.Range("M21").Value = 15.86595001789593
$ws.Range("O21").Value = 31.09474150035393
$ws.Range("B22").Value = 14.9560277527909
$ws.Range("C22").Value = 10.03287451262167
$ws.Range("D22").Value = 8.348901698731519
$ws.Range("F22").Value = 40.50985662234763
$ws.Range("G22").Value = 3.701294578651262
$ws.Range("J22").Value = 11.1894866794701
$ws.Range("K22").Value = 10.72493728949309
$ws.Range("L22").Value = 11.26971523781574
$ws.Range("M22").Value = 15.89520802291011
$ws.Range("O22").Value = 31.04664175542575
$ws.Range("B23").Value = 14.87258447287483
$ws.Range("C23").Value = 10.03140780738581
$ws.Range("D23").Value = 8.341491939098836
$ws.Range("F23").Value = 40.53074626048272
$ws.Range("G23").Value = 3.702011723920263
$ws.Range("J23").Value = 11.19748027412674
$ws.Range("K23").Value = 10.66235854232338
$ws.Range("L23").Value = 11.26975866495666
$ws.Range("M23").Value = 15.87944516569533
$ws.Range("O23").Value = 31.07189834259014
$ws.Range("B24").Value = 14.55619598235048
$ws.Range("C24").Value = 10.02625498260762
$ws.Range("D24").Value = 8.315150139759352
$ws.Range("F24").Value = 40.62281970807989
$ws.Range("G24").Value = 3.704834835637421
$ws.Range("J24").Value = 11.22905394459554
$ws.Range("K24").Value = 10.42438812608075
$ws.Range("L24").Value = 11.27198227132336
$ws.Range("M24").Value = 15.82301177362193
$ws.Range("O24").Value = 31.17666922560824
$ws.Range("B25").Value = 14.21588665048513
$ws.Range("C25").Value = 10.02156989106384
$ws.Range("D25").Value = 8.29050067967755
$ws.Range("F25").Value = 40.74978283361654
$ws.Range("G25").Value = 3.708110281947549
$ws.Range("J25").Value = 11.26590258330987
$ws.Range("K25").Value = 10.16693098075878
$ws.Range("L25").Value = 11.27591688379726
$ws.Range("M25").Value = 15.76929854198071
$ws.Range("O25").Value = 31.2619101836838
